$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.562.91"
$ws.Range("E2").Value = "  -2.54%  "

$ws.Range("D3").Value = "1.752.88"
$ws.Range("E3").Value = "  -3.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4467"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3626"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.034"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.166"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.04%  "

$ws.Range("D16").Value = "1.750.98"
$ws.Range("E16").Value = "  -3.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06380"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.861"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.12%  "

$ws.Range("D23").Value = "27.601.74"
$ws.Range("E23").Value = "  -2.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.102"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "

$ws.Range("D28").Value = "1.952.52"
$ws.Range("E28").Value = "  -3.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.078"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.671"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09015"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.543"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02314"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2092"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6359"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.64%  "

# Row 39: InternetComputer(DFINITY) -> Hedera
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05983"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.35%  "

# Row 40: Hedera -> InternetComputer(DFINITY)
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.964"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.197"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.395"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.769"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5892"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.714"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.955"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.159"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06875"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.74%  "

